$wb = $excel.ActiveWorkbook

# ---- RQ2.1 ----
$ws = $wb.Worksheets.Item("RQ2.1")
$ws.Activate()
$ws.Range("F21").Value = "Device, Data, and Service Authorisation; Interaction and incident record"
$ws.Range("G21").Value = "Fog"
$ws.Range("F27").Value = "Digital Twins of IoT devices"
$ws.Range("G27").Value = "Cloud"
$ws.Range("F34").Value = "Device and Data Integrity validator"
$ws.Range("G34").Value = "Cloud"
$ws.Range("F41").Value = "Communication channel; Device, Data, and Service Authorisation; Device Authentication"
$ws.Range("G41").Value = "Cloud"
$ws.Range("F45").Value = "Sensor Data Storage; Trust management system"
$ws.Range("G45").Value = "Cloud"
$ws.Range("F61").Value = "Business process orchestrator"
$ws.Range("G61").Value = "Fog"
$ws.Range("F66").Value = "Sensor Data Storage"
$ws.Range("G66").Value = "Cloud"
$ws.Range("F67").Value = "Sensor Data Storage; Publish-Subscribe middleware"
$ws.Range("G67").Value = "Cloud"
$ws.Range("F79").Value = "Sensor Data Storage; Business process orchestrator"
$ws.Range("G79").Value = "Cloud"
$ws.Range("F88").Value = "Service registry and matchmaker; Business process orchestrator"
$ws.Range("G88").Value = "Cloud"
$ws.Range("F89").Value = "Device, Data, and Service Authorisation"
$ws.Range("G89").Value = "Cloud"

# ---- RQ2.2 ----
$ws = $wb.Worksheets.Item("RQ2.2")
$ws.Activate()
$ws.Range("F21").Value = "Device interaction records"
$ws.Range("G21").Value = "N/A"
$ws.Range("H21").Value = "Authorisation mechanism"
$ws.Range("I21").Value = "N/A"
$ws.Range("J21").Value = "Trail prediction"
$ws.Range("F27").Value = "Device interaction records"
$ws.Range("G27").Value = "N/A"
$ws.Range("H27").Value = "Digital-twin of devices"
$ws.Range("I27").Value = "N/A"
$ws.Range("J27").Value = "N/A"
$ws.Range("F34").Value = "Commands to devices"
$ws.Range("G34").Value = "N/A"
$ws.Range("H34").Value = "Command integrity check"
$ws.Range("I34").Value = "N/A"
$ws.Range("J34").Value = "N/A"
$ws.Range("F41").Value = "Device interaction records"
$ws.Range("G41").Value = "N/A"
$ws.Range("H41").Value = "Authorisation mechanism"
$ws.Range("I41").Value = "Device Authentication"
$ws.Range("J41").Value = "N/A"
$ws.Range("F45").Value = "Sensor Readings"
$ws.Range("G45").Value = "N/A"
$ws.Range("H45").Value = "Authorisation mechanism; Contract between resource providers and consumers"
$ws.Range("I45").Value = "N/A"
$ws.Range("J45").Value = "Reputation score calculation"
$ws.Range("F61").Value = "Sensor Readings"
$ws.Range("G61").Value = "N/A"
$ws.Range("H61").Value = "Business process"
$ws.Range("I61").Value = "N/A"
$ws.Range("J61").Value = "N/A"
$ws.Range("F66").Value = "Sensor Readings"
$ws.Range("G66").Value = "N/A"
$ws.Range("H66").Value = "N/A"
$ws.Range("I66").Value = "N/A"
$ws.Range("J66").Value = "N/A"
$ws.Range("F67").Value = "Sensor Reading Hashes"
$ws.Range("G67").Value = "N/A"
$ws.Range("H67").Value = "Publish-subscribe protocol"
$ws.Range("I67").Value = "N/A"
$ws.Range("J67").Value = "N/A"
$ws.Range("F79").Value = "Sensor Readings"
$ws.Range("G79").Value = "N/A"
$ws.Range("H79").Value = "N/A"
$ws.Range("I79").Value = "N/A"
$ws.Range("J79").Value = "N/A"
$ws.Range("F88").Value = "Resource exchange records; Service interaction records"
$ws.Range("G88").Value = "N/A"
$ws.Range("H88").Value = "Contract between resource providers and consumers; Service Matchmaking mechanism"
$ws.Range("I88").Value = "N/A"
$ws.Range("J88").Value = "N/A"
$ws.Range("F89").Value = "Authorisation requests and responses"
$ws.Range("G89").Value = "N/A"
$ws.Range("H89").Value = "Authorisation mechanism"
$ws.Range("I89").Value = "N/A"
$ws.Range("J89").Value = "N/A"

# ---- RQ2.3 ----
$ws = $wb.Worksheets.Item("RQ2.3")
$ws.Activate()
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = "blockchain"
$ws.Range("H21").Value = "UTXO"
$ws.Range("I21").Value = "installed"
$ws.Range("J21").Value = "PBFT"
$ws.Range("K21").Value = "Private"
$ws.Range("L21").Value = "local CA"
$ws.Range("M21").Value = "Hyperledger Fabric"
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = "blockchain"
$ws.Range("H27").Value = "account"
$ws.Range("I27").Value = "on-chain"
$ws.Range("J27").Value = "Proof-of-work"
$ws.Range("K27").Value = "public"
$ws.Range("L27").Value = "N/A"
$ws.Range("M27").Value = "Ethereum"
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = "blockchain"
$ws.Range("H34").Value = "account"
$ws.Range("I34").Value = "N/A"
$ws.Range("J34").Value = "Proof-of-work"
$ws.Range("K34").Value = "Private"
$ws.Range("L34").Value = "N/A"
$ws.Range("M34").Value = "Ethereum"
$ws.Range("F41").Value = 1
$ws.Range("G41").Value = "Blockchain"
$ws.Range("H41").Value = "account"
$ws.Range("I41").Value = "on-chain"
$ws.Range("J41").Value = "Proof-of-work"
$ws.Range("K41").Value = "Public"
$ws.Range("L41").Value = "N/A"
$ws.Range("M41").Value = "Ethereum"
$ws.Range("F45").Value = 1
$ws.Range("G45").Value = "blockchain"
$ws.Range("H45").Value = "UTXO"
$ws.Range("I45").Value = "installed"
$ws.Range("J45").Value = "Proof-of-work; Proof-of-space"
$ws.Range("K45").Value = "Consortium"
$ws.Range("L45").Value = "N/A"
$ws.Range("M45").Value = "Hyperledger Fabric"
$ws.Range("F61").Value = 1
$ws.Range("G61").Value = "blockchain"
$ws.Range("H61").Value = "account"
$ws.Range("I61").Value = "on-chain"
$ws.Range("J61").Value = "Proof-of-stake"
$ws.Range("K61").Value = "Public"
$ws.Range("L61").Value = "N/A"
$ws.Range("M61").Value = "Ethereum"
$ws.Range("F66").Value = 1
$ws.Range("G66").Value = "blockchain"
$ws.Range("H66").Value = "UTXO"
$ws.Range("I66").Value = "N/A"
$ws.Range("J66").Value = "Proof-of-work"
$ws.Range("K66").Value = "Public"
$ws.Range("L66").Value = "N/A"
$ws.Range("M66").Value = "In-house BC system"
$ws.Range("F67").Value = 1
$ws.Range("G67").Value = "blockchain"
$ws.Range("H67").Value = "account"
$ws.Range("I67").Value = "on-chain"
$ws.Range("J67").Value = "Proof-of-work"
$ws.Range("K67").Value = "Public"
$ws.Range("L67").Value = "N/A"
$ws.Range("M67").Value = "Ethereum"
$ws.Range("F79").Value = 4
$ws.Range("G79").Value = "blockchain"
$ws.Range("H79").Value = "UTXO"
$ws.Range("I79").Value = "N/A"
$ws.Range("J79").Value = "Proof-of-work"
$ws.Range("K79").Value = "Public"
$ws.Range("L79").Value = "N/A"
$ws.Range("M79").Value = "In-house BC system"
$ws.Range("F88").Value = 1
$ws.Range("G88").Value = "blockchain"
$ws.Range("H88").Value = "account"
$ws.Range("I88").Value = "on-chain"
$ws.Range("J88").Value = "Proof-of-work"
$ws.Range("K88").Value = "public"
$ws.Range("L88").Value = "N/A"
$ws.Range("M88").Value = "Ethereum"
$ws.Range("F89").Value = 1
$ws.Range("G89").Value = "blockchain"
$ws.Range("H89").Value = "account"
$ws.Range("I89").Value = "on-chain"
$ws.Range("J89").Value = "Proof-of-work"
$ws.Range("K89").Value = "Public"
$ws.Range("L89").Value = "N/A"
$ws.Range("M89").Value = "Ethereum"

# ---- RQ3 ----
$ws = $wb.Worksheets.Item("RQ3")
$ws.Activate()
$ws.Range("F21").Value = "N/A"
$ws.Range("G21").Value = "N/A"
$ws.Range("F27").Value = "N/A"
$ws.Range("G27").Value = "N/A"
$ws.Range("F34").Value = "N/A"
$ws.Range("G34").Value = "N/A"
$ws.Range("F41").Value = "N/A"
$ws.Range("G41").Value = "N/A"
$ws.Range("F45").Value = "N/A"
$ws.Range("G45").Value = "N/A"
$ws.Range("F61").Value = "IoT injects too many transactions at too high rate into blockchains"
$ws.Range("G61").Value = "Proof-of-stake"
$ws.Range("F66").Value = "N/A"
$ws.Range("G66").Value = "N/A"
$ws.Range("F67").Value = "N/A"
$ws.Range("G67").Value = "N/A"
$ws.Range("F79").Value = "Blockchain PoW consumes too much energy for IoT use cases."
$ws.Range("G79").Value = "Selective Proof-of-work"
$ws.Range("F88").Value = "N/A"
$ws.Range("G88").Value = "N/A"
$ws.Range("F89").Value = "N/A"
$ws.Range("G89").Value = "N/A"

# ---- Row height adjustments ----
$ws = $wb.Worksheets.Item("RQ2.2")
$ws.Rows.Item(45).RowHeight = 80
$ws.Rows.Item(88).RowHeight = 80

# ---- View state: scroll position + selection ----
$ws = $wb.Worksheets.Item("RQ2.1")
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 40
$ws.Range("G46").Select()

$ws = $wb.Worksheets.Item("RQ2.2")
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 38
$ws.Range("I46").Select()

$ws = $wb.Worksheets.Item("RQ2.3")
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 40
$ws.Range("M46").Select()

$ws = $wb.Worksheets.Item("RQ3")
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 40
$ws.Range("F42").Select()
